# Applies the "Specific Baskets B6 and B7 for semester 7" commit:
#  - Updates the Section_A / Section_B weekly timetable grids so every
#    scheduled slot is tagged with the classroom "[C405]" and several
#    slots are reshuffled to reflect the new room assignment.
#  - Adds two new reporting worksheets: Semester_Rules and
#    Classroom_Utilization.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Section_A timetable updates
# ---------------------------------------------------------------------
$secA = $wb.Worksheets.Item("Section_A")

$secA.Range("B2").Value = "ELECTIVE_B1 [C405]"
$secA.Range("C2").Value = "CS161 [C405]"
$secA.Range("D2").Value = "ELECTIVE_B1 [C405]"
$secA.Range("E2").Value = "EC161 [C405]"

$secA.Range("B3").Value = "CS161 [C405]"
$secA.Range("C3").Value = "Free"
$secA.Range("D3").Value = "DS161 [C405]"
$secA.Range("E3").Value = "MA161 [C405]"
$secA.Range("F3").Value = "MA162 [C405]"

$secA.Range("B5").Value = "MA162 [C405]"
$secA.Range("C5").Value = "Free"
$secA.Range("E5").Value = "Free"

$secA.Range("B6").Value = "DS161 (Tutorial) [C405]"
$secA.Range("C6").Value = "MA161 (Tutorial) [C405]"
$secA.Range("D6").Value = "EC161 (Tutorial) [C405]"
$secA.Range("E6").Value = "Free"
$secA.Range("F6").Value = "ELECTIVE_B1 (Tutorial) [C405]"

$secA.Range("B7").Value = "EC161 [C405]"
$secA.Range("C7").Value = "DS161 [C405]"
$secA.Range("D7").Value = "Free"
$secA.Range("E7").Value = "Free"
$secA.Range("F7").Value = "MA161 [C405]"

$secA.Range("C8").Value = "MA162 (Tutorial) [C405]"
$secA.Range("D8").Value = "CS161 (Tutorial) [C405]"
$secA.Range("E8").Value = "Free"

# ---------------------------------------------------------------------
# 2. Section_B timetable updates
# ---------------------------------------------------------------------
$secB = $wb.Worksheets.Item("Section_B")

$secB.Range("B2").Value = "ELECTIVE_B1 [C405]"
$secB.Range("C2").Value = "EC161 [C405]"
$secB.Range("D2").Value = "ELECTIVE_B1 [C405]"
$secB.Range("E2").Value = "Free"
$secB.Range("F2").Value = "DS161 [C405]"

$secB.Range("B3").Value = "MA162 [C405]"
$secB.Range("C3").Value = "Free"
$secB.Range("D3").Value = "MA162 [C405]"
$secB.Range("E3").Value = "Free"
$secB.Range("F3").Value = "EC161 [C405]"

$secB.Range("B5").Value = "CS161 [C405]"
$secB.Range("C5").Value = "MA161 [C405]"
$secB.Range("D5").Value = "Free"
$secB.Range("E5").Value = "DS161 [C405]"
$secB.Range("F5").Value = "Free"

$secB.Range("C6").Value = "Free"
$secB.Range("D6").Value = "DS161 (Tutorial) [C405]"
$secB.Range("F6").Value = "ELECTIVE_B1 (Tutorial) [C405]"

$secB.Range("C7").Value = "Free"
$secB.Range("D7").Value = "MA161 [C405]"
$secB.Range("E7").Value = "CS161 [C405]"

$secB.Range("B8").Value = "MA161 (Tutorial) [C405]"
$secB.Range("C8").Value = "Free"
$secB.Range("D8").Value = "EC161 (Tutorial) [C405]"
$secB.Range("E8").Value = "MA162 (Tutorial) [C405]"
$secB.Range("F8").Value = "CS161 (Tutorial) [C405]"

# ---------------------------------------------------------------------
# 3. New sheet: Semester_Rules
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$semRules = $wb.Worksheets.Add($null, $lastSheet)
$semRules.Name = "Semester_Rules"

$srHeaders = @("Semester", "Rule", "Exclusion", "Reason", "Scheduled Baskets", "Status")
for ($c = 1; $c -le $srHeaders.Length; $c++) {
    $cell = $semRules.Cells.Item(1, $c)
    $cell.Value = $srHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$semRules.Cells.Item(2, 1).Value = "Semester 1"
$semRules.Cells.Item(2, 2).Value = "Schedule all elective baskets"
$semRules.Cells.Item(2, 3).Value = "None"
$semRules.Cells.Item(2, 4).Value = "No specific restrictions for this semester"
$semRules.Cells.Item(2, 5).Value = "ELECTIVE_B1"
$semRules.Cells.Item(2, 6).Value = [char]0x2705 + " Applied"

# ---------------------------------------------------------------------
# 4. New sheet: Classroom_Utilization
# ---------------------------------------------------------------------
$classUtil = $wb.Worksheets.Add($null, $semRules)
$classUtil.Name = "Classroom_Utilization"

$cuHeaders = @("Room Number", "Type", "Capacity", "Weekly Hours (Timetable)", "Daily Avg Hours (Timetable)", "Exam Sessions", "Utilization Rate (%)", "Facilities")
for ($c = 1; $c -le $cuHeaders.Length; $c++) {
    $cell = $classUtil.Cells.Item(1, $c)
    $cell.Value = $cuHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$cuData = @(
    @("C001", "Recreation", "nil", 0, 0, 0, 0, ""),
    @("C002", "large classroom", "116", 0, 0, 0, 0, "Projector"),
    @("C003", "large classroom", "135", 0, 0, 0, 0, "Projector"),
    @("C004", "Auditorium", "240", 0, 0, 0, 0, "Projector"),
    @("C101", "classroom", "96", 0, 0, 0, 0, "Projector"),
    @("C102", "classroom", "96", 0, 0, 0, 0, "Projector"),
    @("C103", "library", "nil", 0, 0, 0, 0, "Computers"),
    @("C104", "classroom", "96", 0, 0, 0, 0, "Projector"),
    @("L105", "Hardware Lab", "40", 0, 0, 0, 0, "Hardware Equipment"),
    @("L106", "Software Lab", "40", 0, 0, 0, 0, "Computers"),
    @("L107", "Software Lab", "40", 0, 0, 0, 0, "Computers"),
    @("C201", "classroom", "96", 0, 0, 0, 0, "Projector"),
    @("C202", "classroom", "96", 0, 0, 0, 0, "Projector"),
    @("C203", "classroom", "96", 0, 0, 0, 0, "Projector"),
    @("C204", "classroom", "96", 0, 0, 0, 0, "Projector"),
    @("C205", "classroom", "96", 0, 0, 0, 0, "Projector"),
    @("L206", "Hardware Lab", "40", 0, 0, 0, 0, "Hardware Equipment"),
    @("L207", "Software Lab", "40", 0, 0, 0, 0, "Computers"),
    @("L208", "Software Lab", "40", 0, 0, 0, 0, "Computers"),
    @("C301", "Physics Lab", "40", 0, 0, 0, 0, "Projector"),
    @("C302", "classroom", "96", 0, 0, 0, 0, "Projector"),
    @("C303", "classroom", "96", 0, 0, 0, 0, "Projector"),
    @("C304", "classroom", "96", 0, 0, 0, 0, "Projector"),
    @("C305", "classroom", "96", 0, 0, 0, 0, "Projector"),
    @("L306", "classroom", "96", 0, 0, 0, 0, "Computers"),
    @("L307", "Research Scholar Lab", "40", 0, 0, 0, 0, "Computers"),
    @("L308", "Research Scholar Lab", "40", 0, 0, 0, 0, "Computers"),
    @("C401", "classroom", "96", 0, 0, 0, 0, "Projector"),
    @("C402", "classroom", "96", 0, 0, 0, 0, "Projector"),
    @("C403", "classroom", "78", 0, 0, 0, 0, "Projector"),
    @("C404", "classroom", "78", 0, 0, 0, 0, "Projector"),
    @("C405", "classroom", "78", 48, 9.6, 0, 100, "Projector"),
    @("L406", "classroom", "78", 0, 0, 0, 0, "Computers"),
    @("L407", "classroom", "78", 0, 0, 0, 0, "Computers"),
    @("L408", "classroom without projector", "78", 0, 0, 0, 0, "Computers")
)

$r = 2
foreach ($row in $cuData) {
    $classUtil.Cells.Item($r, 1).Value = $row[0]
    $classUtil.Cells.Item($r, 2).Value = $row[1]
    $classUtil.Cells.Item($r, 3).NumberFormat = "@"
    $classUtil.Cells.Item($r, 3).Value = $row[2]
    $classUtil.Cells.Item($r, 4).Value = $row[3]
    $classUtil.Cells.Item($r, 5).Value = $row[4]
    $classUtil.Cells.Item($r, 6).Value = $row[5]
    $classUtil.Cells.Item($r, 7).Value = $row[6]
    $classUtil.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
